$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty "autogluon" row (row 3) results.
$ws.Range("B3").Value = "0.421 (0.383 ± 0.019)"
$ws.Range("C3").Value = "00:01:54 (00:02:06 ± 00:00:06)"
$ws.Range("D3").Value = "00:00:00 (00:00:00 ± 00:00:00)"
$ws.Range("E3").Value = "[]"
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "31"
$ws.Range("F3").Style = "Normal"

# Repair the mis-decoded "±" (was double-encoded as "Â±") in the already
# populated rows.
$ws.Range("B4").Value = "0.639 (0.567 ± 0.037)"
$ws.Range("C4").Value = "00:01:48 (00:02:20 ± 00:00:32)"
$ws.Range("D4").Value = "00:00:00 (00:00:00 ± 00:00:00)"

$ws.Range("B6").Value = "0.658 (0.625 ± 0.020)"
$ws.Range("C6").Value = "00:04:57 (00:05:01 ± 00:00:01)"
$ws.Range("D6").Value = "00:00:00 (00:00:02 ± 00:00:01)"

$ws.Range("B8").Value = "0.538 (0.489 ± 0.027)"
$ws.Range("C8").Value = "00:04:56 (00:05:56 ± 00:00:34)"
$ws.Range("D8").Value = "00:00:00 (00:00:00 ± 00:00:00)"
